$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New exception-handling / payment-status lookup rows appended below the
# existing settings table (rows 9-12), each with Name / Value / Description.
# Names and values are written first, and the shared Description text last,
# so new shared-string entries are created in the same order as the target.
$ws.Range("A9").Value = "Statuspaid"
$ws.Range("B9").Value = "Paid"

$ws.Range("A10").Value = "StatusUnPaid"
$ws.Range("B10").Value = "Unpaid"

$ws.Range("A11").Value = "StatusPartiallyPaid"
$ws.Range("B11").Value = "Partially Paid"

$ws.Range("A12").Value = "StatusOverDue"
$ws.Range("B12").Value = "Overdue"

$ws.Range("C9").Value = "Payment status matching with status dropdown in ERP"
$ws.Range("C10").Value = "Payment status matching with status dropdown in ERP"
$ws.Range("C11").Value = "Payment status matching with status dropdown in ERP"
$ws.Range("C12").Value = "Payment status matching with status dropdown in ERP"

# Move / record the active selection as it was left after the edit.
[void]$ws.Range("B17").Select()
